$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("params")

# Delete column D ("policy.EL") - shifts everything left
$ws1.Columns.Item(4).Delete()

# Rename old column C header (now holds "policy.SR") to "policy"
$ws1.Range("C4").Value = "policy"

# Copy row 5 down into rows 8, 9, 10 as base data for new scenarios
$ws1.Range("A5:AJ5").Copy($ws1.Range("A8:AJ8"))
$ws1.Range("A5:AJ5").Copy($ws1.Range("A9:AJ9"))
$ws1.Range("A5:AJ5").Copy($ws1.Range("A10:AJ10"))

# the source row has no values in D or AI:AJ - clear the copies so no
# stray blank cells get materialized there
$ws1.Range("D8:D10").ClearContents()
$ws1.Range("AI8:AJ10").ClearContents()

$ws1.Range("A8").Value = "RS1_base"
$ws1.Range("A9").Value = "RS2_base"
$ws1.Range("A10").Value = "RS3_base"

$ws1.Range("F8:F10").Value = $false

# Extend the data validation ranges to cover the shifted columns and the
# newly added rows 8:10
$ws1.Range("F5:H7").Validation.Delete()
$ws1.Range("Z5:Z7").Validation.Delete()
$ws1.Range("F5:H10").Validation.Add(3, 1, 1, '"TRUE, FALSE"')
$ws1.Range("Z5:Z10").Validation.Add(3, 1, 1, '"simple, internal"')

$ws1.Activate()
$ws1.Range("E25").Select()

$ws3 = $wb.Worksheets.Item("returns")
$ws3.Activate()
$ws3.Range("D20").Select()

$ws1.Activate()
